# Update cryptos list prices and 1h volume percentages
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "27.919.86"
$ws.Range("E2").Value = "  +1.46%  "
$ws.Range("D3").Value = "1.648.21"
$ws.Range("E3").Value = "  +1.78%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "213.82"
$ws.Range("E5").Value = "  +1.29%  "
$ws.Range("E6").Value = "  +0.53%  "
$ws.Range("E7").Value = "  -0.07%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "23.54"
$ws.Range("E8").Value = "  +3.66%  "
$ws.Range("E9").Value = "  +1.02%  "
$ws.Range("E10").Value = "  +0.21%  "
$ws.Range("E11").Value = "  -1.48%  "
$ws.Range("E12").Value = "  +1.74%  "
$ws.Range("D13").Value = "1.646.22"
$ws.Range("E13").Value = "  +1.71%  "
$ws.Range("E14").Value = "  +1.51%  "
$ws.Range("E15").Value = "  +2.25%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "65.65"
$ws.Range("E16").Value = "  +1.02%  "
$ws.Range("D17").Value = "27.920.21"
$ws.Range("E17").Value = "  +1.45%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "232.34"
$ws.Range("E18").Value = "  +1.09%  "
$ws.Range("E19").Value = "  +1.87%  "
$ws.Range("E20").Value = "  +0.51%  "
$ws.Range("E21").Value = "  -0.07%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "10.74"
$ws.Range("E22").Value = "  +5.80%  "
$ws.Range("E23").Value = "  +2.48%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "152.33"
$ws.Range("E25").Value = "  +1.69%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "6.90"
$ws.Range("E26").Value = "  +1.33%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "15.72"
$ws.Range("E27").Value = "  +1.07%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "0.111"
$ws.Range("E28").Value = "  +0.27%  "
$ws.Range("E29").Value = "  -0.07%  "
$ws.Range("E30").Value = "  +1.73%  "
$ws.Range("E31").Value = "  +0.58%  "
$ws.Range("E32").Value = "  +2.63%  "
$ws.Range("D33").Value = "1.451.67"
$ws.Range("E33").Value = "  +0.62%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "3.11"
$ws.Range("E34").Value = "  +1.58%  "
$ws.Range("E35").Value = "  +1.80%  "
$ws.Range("E36").Value = "  -0.55%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.891"
$ws.Range("E37").Value = "  +3.36%  "
$ws.Range("E38").Value = "  +0.74%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.561"
$ws.Range("E39").Value = "  -0.05%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.921"
$ws.Range("E40").Value = "  -1.54%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "69.34"
$ws.Range("E41").Value = "  +0.54%  "
$ws.Range("E42").Value = "  +2.39%  "
$ws.Range("E43").Value = "  -0.07%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "2.46"
$ws.Range("E44").Value = "  +0.41%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "2.24"
$ws.Range("E45").Value = "  +1.21%  "
$ws.Range("E46").Value = "  -0.64%  "
$ws.Range("E47").Value = "  +5.42%  "
$ws.Range("D48").Value = "1.789.61"
$ws.Range("E48").Value = "  +1.67%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "88.79"
$ws.Range("E49").Value = "  +2.99%  "
$ws.Range("E50").Value = "  +0.37%  "
$ws.Range("E51").Value = "  +1.05%  "
